# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Ravana_Profits workbook
# (columns H-N are cached market-board data with no formulas; this mirrors a data refresh)

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2513.7273
$ws.Range("I19").Value = 2700.2856
$ws.Range("K19").Value = 2700.2856
$ws.Range("M19").Value = -2525.2856
$ws.Range("H132").Value = 1126
$ws.Range("I132").Value = 1126
$ws.Range("K132").Value = 3378
$ws.Range("M132").Value = -848
$ws.Range("H137").Value = 1823.129
$ws.Range("I137").Value = 1409.95
$ws.Range("J137").Value = 2574.3635
$ws.Range("K137").Value = 4229.85
$ws.Range("L137").Value = 7723.0905
$ws.Range("M137").Value = -1679.85
$ws.Range("N137").Value = -12823.0905
$ws.Range("H138").Value = 3532.3167
$ws.Range("I138").Value = 3142
$ws.Range("J138").Value = 3583.868
$ws.Range("K138").Value = 9426
$ws.Range("L138").Value = 10751.604
$ws.Range("M138").Value = -4286
$ws.Range("N138").Value = -21031.604

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5380.203
$ws.Range("I32").Value = 4141.2373
$ws.Range("K32").Value = 4141.2373
$ws.Range("M32").Value = -3854.2373
$ws.Range("H63").Value = 1198.7142
$ws.Range("I63").Value = 1064.8334
$ws.Range("K63").Value = 1064.8334
$ws.Range("M63").Value = -378.8334
$ws.Range("H66").Value = 1198.7142
$ws.Range("I66").Value = 1064.8334
$ws.Range("K66").Value = 5324.166999999999
$ws.Range("M66").Value = -1892.166999999999
$ws.Range("H102").Value = 4166.6665
$ws.Range("I102").Value = 4166.6665
$ws.Range("K102").Value = 4166.6665
$ws.Range("M102").Value = -2544.6665
$ws.Range("H122").Value = 5781.5
$ws.Range("I122").Value = 6137.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 18413.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -15963.4
$ws.Range("N122").Value = -16900

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 44861
$ws.Range("J55").Value = 44861
$ws.Range("L55").Value = 44861
$ws.Range("N55").Value = -45407
$ws.Range("H64").Value = 1326.6666
$ws.Range("J64").Value = 1490
$ws.Range("L64").Value = 1490
$ws.Range("N64").Value = -1940
$ws.Range("H67").Value = 1326.6666
$ws.Range("J67").Value = 1490
$ws.Range("L67").Value = 1490
$ws.Range("N67").Value = -3050
$ws.Range("H99").Value = 731.1667
$ws.Range("I99").Value = 731.1667
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 731.1667
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 766.8333
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H134").Value = 2776.8125
$ws.Range("I134").Value = 2735.7
$ws.Range("J134").Value = 2845.3333
$ws.Range("K134").Value = 8207.099999999999
$ws.Range("L134").Value = 8535.999899999999
$ws.Range("M134").Value = -5672.099999999999
$ws.Range("N134").Value = -13605.9999

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 906.25
$ws.Range("I7").Value = 908.3333
$ws.Range("J7").Value = 900
$ws.Range("K7").Value = 908.3333
$ws.Range("L7").Value = 900
$ws.Range("M7").Value = -795.3333
$ws.Range("N7").Value = -1126
$ws.Range("H31").Value = 1824.375
$ws.Range("I31").Value = 1561.75
$ws.Range("J31").Value = 2087
$ws.Range("K31").Value = 1561.75
$ws.Range("L31").Value = 2087
$ws.Range("M31").Value = -1266.75
$ws.Range("N31").Value = -2677
$ws.Range("H34").Value = 1824.375
$ws.Range("I34").Value = 1561.75
$ws.Range("J34").Value = 2087
$ws.Range("K34").Value = 1561.75
$ws.Range("L34").Value = 2087
$ws.Range("M34").Value = -1359.75
$ws.Range("N34").Value = -2491

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 244.75
$ws.Range("J2").Value = 490
$ws.Range("L2").Value = 2940
$ws.Range("N2").Value = -3166
$ws.Range("H69").Value = 1332.6666
$ws.Range("I69").Value = 1224.5
$ws.Range("K69").Value = 3673.5
$ws.Range("M69").Value = -2862.5
$ws.Range("H72").Value = 1332.6666
$ws.Range("I72").Value = 1224.5
$ws.Range("K72").Value = 11020.5
$ws.Range("M72").Value = -6964.5
$ws.Range("H75").Value = 1300
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1300
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H82").Value = 19000
$ws.Range("I82").Value = 19000
$ws.Range("K82").Value = 57000
$ws.Range("M82").Value = -56594
$ws.Range("H85").Value = 19000
$ws.Range("I85").Value = 19000
$ws.Range("K85").Value = 57000
$ws.Range("M85").Value = -55596
$ws.Range("H138").Value = 7942.6665
$ws.Range("I138").Value = 5885.6665
$ws.Range("J138").Value = 9999.666999999999
$ws.Range("K138").Value = 17656.9995
$ws.Range("L138").Value = 29999.001
$ws.Range("M138").Value = -12516.9995
$ws.Range("N138").Value = -40279.001

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108.625
$ws.Range("I2").Value = 103.166664
$ws.Range("K2").Value = 103.166664
$ws.Range("M2").Value = 9.833336000000003
$ws.Range("H97").Value = 398.37036
$ws.Range("I97").Value = 484
$ws.Range("J97").Value = 98.666664
$ws.Range("K97").Value = 484
$ws.Range("L97").Value = 98.666664
$ws.Range("M97").Value = 12
$ws.Range("N97").Value = -1090.666664
$ws.Range("H126").Value = 999
$ws.Range("I126").Value = 998
$ws.Range("K126").Value = 2994
$ws.Range("M126").Value = -524
$ws.Range("H132").Value = 4268.769
$ws.Range("I132").Value = 4166
$ws.Range("J132").Value = 4299.6
$ws.Range("K132").Value = 12498
$ws.Range("L132").Value = 12898.8
$ws.Range("M132").Value = -9968
$ws.Range("N132").Value = -17958.8
$ws.Range("H134").Value = 110326
$ws.Range("J134").Value = 110326
$ws.Range("L134").Value = 330978
$ws.Range("N134").Value = -336048

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2004
$ws.Range("I40").Value = 2004
$ws.Range("K40").Value = 2004
$ws.Range("M40").Value = -1868
$ws.Range("H93").Value = 3396.5
$ws.Range("I93").Value = 3396.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3396.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2148.5
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 9966.5
$ws.Range("I122").Value = 7899.5
$ws.Range("K122").Value = 23698.5
$ws.Range("M122").Value = -21248.5
$ws.Range("H137").Value = 90450
$ws.Range("J137").Value = 90450
$ws.Range("L137").Value = 90450
$ws.Range("N137").Value = -100650

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1492.5
$ws.Range("I96").Value = 1492.5
$ws.Range("K96").Value = 1492.5
$ws.Range("M96").Value = -119.5
$ws.Range("H136").Value = 2719.25
$ws.Range("I136").Value = 2546.5264
$ws.Range("K136").Value = 7639.5792
$ws.Range("M136").Value = -5089.5792
